$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pool")
$ws.Activate()

$ws.Range("AP1").Value = "VL"
$ws.Range("AQ1").Value = "SL"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 42).Value = "15"
    $ws.Cells.Item($r, 43).Value = "15"
}

[void]$ws.Range("AG14").Select()
